$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.897.50'
$ws.Range("D3").Value = '3.313.83'
$ws.Range("E3").Value = '  +6.02%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.47'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.25'
$ws.Range("E6").Value = '  +5.31%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.312.16'
$ws.Range("E8").Value = '  +6.33%  '
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("E10").Value = '  +2.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.50'
$ws.Range("E11").Value = '  +5.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.475'
$ws.Range("E12").Value = '  +3.80%  '
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.86'
$ws.Range("E14").Value = '  +2.21%  '
$ws.Range("D15").Value = '3.856.74'
$ws.Range("E15").Value = '  +5.98%  '
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").Value = '3.311.99'
$ws.Range("E17").Value = '  +5.99%  '
$ws.Range("D18").Value = '63.974.92'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.64'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("E22").Value = '  +6.35%  '
$ws.Range("E23").Value = '  +4.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.62'
$ws.Range("E24").Value = '  +4.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.77'
$ws.Range("E25").Value = '  -2.62%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  +2.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +1.97%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +4.00%  '
$ws.Range("E31").Value = '  +5.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.05'
$ws.Range("E32").Value = '  +8.90%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("E34").Value = '  +1.65%  '
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  +3.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.33'
$ws.Range("E37").Value = '  +2.55%  '
$ws.Range("D38").Value = '0.0₃0757'
$ws.Range("E38").Value = '  +7.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0399'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '431.90'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '3.037.01'
$ws.Range("E41").Value = '  +5.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.46'
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("E43").Value = '  +4.22%  '
$ws.Range("E44").Value = '  -3.88%  '
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("E46").Value = '  +4.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.52'
$ws.Range("E47").Value = '  +3.11%  '
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.54'
$ws.Range("E50").Value = '  +12.02%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.32'
$ws.Range("E51").Value = '  +1.88%  '
